$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-24 down to 7-25
$ws.Rows(6).Insert()

# Populate the new row 6 with the new weekly price record
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 44757
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = 100112043
$ws.Cells.Item(6, 7).Value = "Pepino dulce"
$ws.Cells.Item(6, 8).Value = "Cultivar XV región"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 150
$ws.Cells.Item(6, 11).Value = 6000
$ws.Cells.Item(6, 12).Value = 6500
$ws.Cells.Item(6, 13).Value = 6250
$ws.Cells.Item(6, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(6, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 16).Value = 625
$ws.Cells.Item(6, 17).Value = 10
$ws.Cells.Item(6, 18).Value = "Hortaliza"
